# TUERCAS price list refresh (fix bug "exceeded requests" in google drive):
#  - bump the sheet date in A1 by one day
#  - refresh the zinc-plated nut prices in column D
#  - keep the merged title range (A1:E1) registered ahead of the A30:D30
#    merge, as in the target workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds a date serial; advance it by one day (45310 -> 45311).
# Value2 is used instead of Value because Value returns a COM variant
# wrapper in this host that does not coerce cleanly in arithmetic.
$ws.Range("A1").Value2 = $ws.Range("A1").Value2 + 1

# Updated "PRECIO ZINC." prices for the altas (TUA-10x) rows
$ws.Range("D23").Value2 = 4848.095
$ws.Range("D24").Value2 = 6683.032
$ws.Range("D25").Value2 = 9503.041999999999
$ws.Range("D26").Value2 = 14177.304
$ws.Range("D27").Value2 = 25573.226
$ws.Range("D28").Value2 = 33106.119

# Updated prices for the bajas (TU-10x) rows
$ws.Range("D36").Value2 = 3136
$ws.Range("D37").Value2 = 4696

# The workbook XML re-orders the merged-cell list so that A1:E1 is
# registered before A30:D30. Re-register the A30:D30 merge so it is
# re-appended after the existing A1:E1 entry, reproducing that order.
$ws.Range("A30:D30").UnMerge()
$ws.Range("A30:D30").Merge()
